$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "'-3.03%"
$ws.Range("D3").Value = "'54.20"
$ws.Range("E3").Value = "'9.67%"
$ws.Range("D4").Value = "'5.095"
$ws.Range("E4").Value = "'-4.20%"
$ws.Range("D5").Value = "'0.07918"
$ws.Range("E6").Value = "'-0.36%"
$ws.Range("D7").Value = "'1.397"
$ws.Range("E7").Value = "'4.28%"
$ws.Range("D8").Value = "'1.663"
$ws.Range("E8").Value = "'1.02%"
$ws.Range("D9").Value = "'0.1234"
$ws.Range("E9").Value = "'-3.90%"
$ws.Range("D10").Value = "'0.2003"
$ws.Range("E10").Value = "'1.69%"
$ws.Range("D11").Value = "'0.04729"
$ws.Range("E11").Value = "'0.65%"
$ws.Range("D12").Value = "'0.09467"
$ws.Range("E12").Value = "'-2.44%"
$ws.Range("D13").Value = "'0.1042"
$ws.Range("E13").Value = "'-0.51%"
$ws.Range("D14").Value = "'0.001279"
$ws.Range("E14").Value = "'-3.07%"
$ws.Range("D15").Value = "'0.005832"
$ws.Range("E15").Value = "'-1.53%"
$ws.Range("D16").Value = "'3.336"
$ws.Range("E16").Value = "'-0.34%"
$ws.Range("D17").Value = "'2.437"
$ws.Range("E17").Value = "'-0.26%"
$ws.Range("E18").Value = "'-1.14%"
$ws.Range("D19").Value = "'8.382"
$ws.Range("E19").Value = "'3.93%"
$ws.Range("D20").Value = "'0.1359"
$ws.Range("E20").Value = "'-1.55%"
$ws.Range("E21").Value = "'-5.31%"
$ws.Range("D22").Value = "'0.04183"
$ws.Range("E22").Value = "'-0.26%"
$ws.Range("E23").Value = "'-2.98%"
$ws.Range("D24").Value = "'0.003988"
$ws.Range("E24").Value = "'-7.70%"
$ws.Range("E25").Value = "'-0.11%"
$ws.Range("E26").Value = "'-0.09%"
$ws.Range("D38").Value = "'0.02629"
$ws.Range("E38").Value = "'-3.43%"
$ws.Range("D39").Value = "'0.05936"
$ws.Range("E39").Value = "'-1.24%"
$ws.Range("D40").Value = "'0.01080"
$ws.Range("E40").Value = "'0.23%"
$ws.Range("E41").Value = "'15.38%"
$ws.Range("D42").Value = "'0.007962"
$ws.Range("E42").Value = "'-0.59%"
$ws.Range("D43").Value = "'0.008210"
$ws.Range("E43").Value = "'3.90%"
$ws.Range("D44").Value = "'0.008275"
$ws.Range("E44").Value = "'5.06%"
$ws.Range("D45").Value = "'0.3444"
$ws.Range("E45").Value = "'-2.09%"
$ws.Range("D46").Value = "'0.00007303"
$ws.Range("E46").Value = "'3.26%"
$ws.Range("E47").Value = "'-0.10%"
$ws.Range("B48").Value = "BOLO"
$ws.Range("C48").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D48").Value = "'0.05547"
$ws.Range("E48").Value = "'0.57%"
$ws.Range("B49").Value = "CoinbaseStockToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D49").Value = "'0.002613"
$ws.Range("E49").Value = "'-34.57%"
$ws.Range("D50").Value = "'0.00002095"
$ws.Range("E50").Value = "'-0.10%"
$ws.Range("E51").Value = "'-0.10%"
